$wb = $excel.ActiveWorkbook
$ws3 = $wb.Worksheets.Item("回收")

# Header row
$ws3.Range("A1").Value = "物品"
$ws3.Range("B1").Value = "个数"
$ws3.Range("C1").Value = "soj"

# Existing two rows (24# / 25#) gain quantity + soj value columns
$ws3.Range("A2").Value = "24#"
$ws3.Range("B2").Value = 6
$ws3.Range("C2").Value = 240

$ws3.Range("A3").Value = "25#"
$ws3.Range("B3").Value = 1
$ws3.Range("C3").Value = 85

# New item rows
$ws3.Range("A4").Value = "eth15ed战枪6s"
$ws3.Range("B4").Value = 1
$ws3.Range("C4").Value = 450

$ws3.Range("A5").Value = "有形猛禽爪"
$ws3.Range("B5").Value = 1
$ws3.Range("C5").Value = 5

$ws3.Range("A6").Value = "女族长标枪 有形"
$ws3.Range("B6").Value = 1
$ws3.Range("C6").Value = 5

$ws3.Range("A7").Value = "洞穴巨魔巢穴盾牌"
$ws3.Range("B7").Value = 1
$ws3.Range("C7").Value = 5

$ws3.Range("A8").Value = "28#"
$ws3.Range("B8").Value = 1
$ws3.Range("C8").Value = 630

# Total row
$ws3.Range("B9").Value = "合计"
$ws3.Range("C9").Formula = "=SUM(C2:C8)"

# Separate note row
$ws3.Range("A11").Value = "chentuhuishou"
$ws3.Range("B11").Value = 123123

# Column A width (~19 characters)
$ws3.Columns.Item(1).ColumnWidth = 128/7

# Page setup: match portrait / letter-size paper used by the other sheets
$ws3.PageSetup.PaperSize = 9
$ws3.PageSetup.Orientation = 1

# Selection + active sheet -- makes "回收" (sheet index 3) the active tab
# and drops tabSelected from "Sheet1" automatically.
$ws3.Range("A1:C9").Select()
$ws3.Activate()
